# Applies the "Updated cryptos list" price/volume refresh (+ a rank-42/43 coin swap)
# to the cryptos worksheet. Values look numeric (e.g. "36.687.36", "43.40") but are
# stored as TEXT in the source data, so each Price (column D) write is done with a
# leading apostrophe to stop Excel's automatic number coercion (which would both change
# the cell type and silently drop meaningful trailing/formatting digits), then the
# quote-prefix style is cleared so the cell is left in the default "Normal" style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''36.687.36'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.60%  '
# Row 3
$ws.Range('D3').Value = '''1.967.60'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.45%  '
# Row 4
$ws.Range('E4').Value = '  +0.06%  '
# Row 5
$ws.Range('D5').Value = '''244.14'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.52%  '
# Row 6
$ws.Range('D6').Value = '''0.616'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.91%  '
# Row 7
$ws.Range('D7').Value = '''58.51'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.69%  '
# Row 8
$ws.Range('E8').Value = '  +0.09%  '
# Row 9
$ws.Range('D9').Value = '''0.373'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.42%  '
# Row 10
$ws.Range('D10').Value = '''0.0808'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.11%  '
# Row 11
$ws.Range('E11').Value = '  -0.07%  '
# Row 12
$ws.Range('D12').Value = '''22.06'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.16%  '
# Row 13
$ws.Range('D13').Value = '''2.257.44'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.85%  '
# Row 14
$ws.Range('E14').Value = '  +0.47%  '
# Row 15
$ws.Range('D15').Value = '''13.71'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.95%  '
# Row 16
$ws.Range('D16').Value = '''5.27'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.82%  '
# Row 17
$ws.Range('D17').Value = '''1.962.21'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.59%  '
# Row 18
$ws.Range('D18').Value = '''36.635.66'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.80%  '
# Row 19
$ws.Range('D19').Value = '''69.68'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.44%  '
# Row 20
$ws.Range('D20').Value = '''0.0₃0860'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.62%  '
# Row 21
$ws.Range('E21').Value = '  +2.02%  '
# Row 22
$ws.Range('D22').Value = '''227.95'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.40%  '
# Row 23
$ws.Range('E23').Value = '  -0.20%  '
# Row 24
$ws.Range('D24').Value = '''2.39'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.43%  '
# Row 25
$ws.Range('E25').Value = '  +1.91%  '
# Row 26
$ws.Range('D26').Value = '''9.39'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +1.10%  '
# Row 27
$ws.Range('D27').Value = '''160.45'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.14%  '
# Row 28
$ws.Range('E28').Value = '  +10.84%  '
# Row 29
$ws.Range('D29').Value = '''19.34'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.18%  '
# Row 30
$ws.Range('E30').Value = '  +1.57%  '
# Row 31
$ws.Range('E31').Value = '  -1.87%  '
# Row 32
$ws.Range('D32').Value = '''4.69'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.68%  '
# Row 33
$ws.Range('E33').Value = '  -1.24%  '
# Row 34
$ws.Range('E34').Value = '  -0.85%  '
# Row 35
$ws.Range('E35').Value = '  -0.07%  '
# Row 36
$ws.Range('D36').Value = '''6.10'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.31%  '
# Row 37
$ws.Range('D37').Value = '''3.39'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +16.93%  '
# Row 38
$ws.Range('D38').Value = '''2.23'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.06%  '
# Row 39
$ws.Range('E39').Value = '  -0.25%  '
# Row 40
$ws.Range('D40').Value = '''0.0998'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.87%  '
# Row 41
$ws.Range('E41').Value = '  +1.01%  '
# Row 42
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').Value = '''1.17'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.49%  '
# Row 43
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').Value = '''0.0211'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.87%  '
# Row 44
$ws.Range('D44').Value = '''16.07'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.71%  '
# Row 45
$ws.Range('D45').Value = '''1.368.67'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.79%  '
# Row 46
$ws.Range('E46').Value = '  +0.66%  '
# Row 47
$ws.Range('D47').Value = '''87.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.02%  '
# Row 48
$ws.Range('D48').Value = '''7.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.19%  '
# Row 49
$ws.Range('E49').Value = '  +0.96%  '
# Row 50
$ws.Range('D50').Value = '''2.148.38'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.86%  '
# Row 51
$ws.Range('D51').Value = '''43.40'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.87%  '
